$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename the worksheet from "Sheet1" to "user_stats"
$ws.Name = "user_stats"

# Re-affirm the E2 value (floating point representation of the same number)
$ws.Range("E2").Value = 17.100000000000001
